# Add a new weekly record at the top of the data (row 26), shifting the
# existing rows (26-105) down by one row to (27-106).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = 5
$ws.Range("B26").Value = 'Macroferia Regional de Talca'
$ws.Range("C26").Value = 'Maule'
$ws.Range("D26").Value = "09/14/2021"
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112017
$ws.Range("G26").Value = 'Apio'
$ws.Range("H26").Value = 'Americana (o)'
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 600
$ws.Range("K26").Value = 7000
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 7000
$ws.Range("N26").Value = '$/docena de matas'
$ws.Range("O26").Value = 'Provincia del Elquí'
$ws.Range("P26").Value = 1167
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = 'Hortaliza'
